# Apply trade #15 (row 16) closing update to the live trading results workbook.
# - Summary sheet: update current capital / total P&L / trade counts / win rate
# - Strategy Status sheet: update MarketMaking strategy row (row 4)
# - All Trades / MarketMaking sheets: append the new closed trade as row 16

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.93        # Current Capital
$summary.Range("B4").Value = -0.07          # Total P&L $
$summary.Range("B6").Value = 15             # Total Trades
$summary.Range("B8").Value = 9              # Losing Trades
$summary.Range("B9").Value = 33.33          # Win Rate %

# ---------------------------------------------------------------------------
# Sheet: Strategy Status (MarketMaking is row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.93           # Capital
$status.Range("D4").Value = 15              # Trades
$status.Range("E4").Value = -0.07           # P&L $
$status.Range("F4").Value = -0.07           # P&L %
$status.Range("G4").Value = 33.33           # Win Rate %

# ---------------------------------------------------------------------------
# Sheets: "All Trades" and "MarketMaking" both get the new trade row (row 16)
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A16").Value = 15

    # Force text number format on the Date column so Excel doesn't coerce the
    # "2026-02-17" string into a date serial number.
    $ws.Range("B16").NumberFormat = "@"
    $ws.Range("B16").Value = "2026-02-17"

    $ws.Range("C16").Value = "07:53:59"
    $ws.Range("D16").Value = "MarketMaking"
    $ws.Range("E16").Value = "UP"
    $ws.Range("F16").Value = 0.76
    $ws.Range("G16").Value = 0.75
    $ws.Range("H16").Value = "CLOSED"
    $ws.Range("I16").Value = -1.3158
    $ws.Range("J16").Value = -0.01
    $ws.Range("K16").Value = 99.93
    $ws.Range("L16").Value = 0
    $ws.Range("M16").Value = 0
    $ws.Range("N16").Value = 0.6
    $ws.Range("O16").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P16").Value = "early_exit"
    $ws.Range("Q16").Value = 0.14
}
